# Add a new review row (row 15) to Sheet1, mirroring the structure of the
# existing rows, with new appid/keyword/email/recovery-email/time/review data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Populate the new row's cell values (this also creates the three new
#    shared-string entries in the same order they appear in the target file:
#    appid, keyword, review text).
$ws.Range("A15").Value = "com.hamxa.shaynachim"
$ws.Range("B15").Value = "bitcoin guide"
$ws.Range("C15").Value = "eligitel@gmail.com"
$ws.Range("D15").Value = "ronenchen27@gmail.com"
$ws.Range("E15").Value = "27/5/2019 15:59"
$ws.Range("F15").Value = "great bitcoin beginners app with great information"

# 2. Turn the email columns into mailto hyperlinks, matching the existing
#    rows' hyperlink pattern.
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:eligitel@gmail.com", "", "", "eligitel@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:ronenchen27@gmail.com", "", "", "ronenchen27@gmail.com") | Out-Null

# 3. Re-apply the formatting (styles) of the row above (row 14) onto the new
#    row, since adding hyperlinks overrides the cell style with the default
#    "Hyperlink" style. This restores the plain style used by every other
#    data row for columns C/D while keeping the rest identical too.
$ws.Range("A14:F14").Copy() | Out-Null
$ws.Range("A15:F15").PasteSpecial(-4122) | Out-Null

# 4. Move the active selection to F15, matching the saved workbook state.
$ws.Range("F15").Select() | Out-Null
